$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 10 new lookup rows (256-265) describing N-limited livestock outputs
# (dairy cattle / beef cattle / goats & sheep) to the "Appendix A" table.
# ---------------------------------------------------------------------------

# Column A ("PyOrator variable") - fill first, top to bottom, to match the
# order new strings were originally authored in.
$ws.Range("A256").Value = "dairy_cat_milk_prod_nlim"
$ws.Range("A257").Value = "dairy_cat_meat_prod_nlim"
$ws.Range("A258").Value = "dairy_cat_manure_prod_nlim"
$ws.Range("A259").Value = "beef_cat_n_excrete_nlim"
$ws.Range("A260").Value = "beef_cat_meat_prod_nlim"
$ws.Range("A261").Value = "beef_cat_manure_prod_nlim"
$ws.Range("A262").Value = "goats_sheep_n_excrete_nlim"
$ws.Range("A263").Value = "goats_sheep_milk_prod_nlim"
$ws.Range("A264").Value = "goats_sheep_meat_prod_nlim"
$ws.Range("A265").Value = "goats_sheep_manure_prod_nlim"

# Column C ("PyOrator display") / Column E ("Definition") - the N-excretion
# rows (259, 262) were documented first, then the remaining display labels,
# then the remaining definitions.
$ws.Range("C259").Value = "N excreted by Beef Cattle (N Limited)"
$ws.Range("C262").Value = "N excreted by goats and sheep (N Limited)"
$ws.Range("E259").Value = "Nitrogen excreted by beef cattle, crop growth limited by N availability"
$ws.Range("E262").Value = "Nitrogen excreted by goats and sheep , crop growth limited by N availability"

$ws.Range("C256").Value = "Milk produced by dairy cattle (N Limited)"
$ws.Range("C257").Value = "Meat produced by dairy cattle (N Limited)"
$ws.Range("C258").Value = "Manure produced by dairy cattle (N Limited)"
$ws.Range("C260").Value = "Meat produced by beef cattle (N Limited)"
$ws.Range("C261").Value = "Manure produced by beef cattle (N Limited)"
$ws.Range("C263").Value = "Milk produced by goats and sheep (N Limited)"

$ws.Range("E256").Value = "Milk produced by dairy cattle, crop growth limited by N availability"
$ws.Range("E263").Value = "Milk produced by goats and sheep, crop growth limited by N availability"
$ws.Range("E257").Value = "Meat produced by dairy cattle, crop growth limited by N availability"
$ws.Range("E260").Value = "Meat produced by beef cattle, crop growth limited by N availability"
$ws.Range("E261").Value = "Manure produced by beef cattle, crop growth limited by N availability"
$ws.Range("E264").Value = "Meat produced by  goats and sheep, crop growth limited by N availability"

$ws.Range("C265").Value = "Manure produced by goats and sheep (N Limited)"
$ws.Range("C264").Value = "Meat produced by goats and sheep (N Limited)"

$ws.Range("E258").Value = "Manure produced by dairy cattle, crop growth limited by N availability"
$ws.Range("E265").Value = "Manure produced by goats and sheep, , crop growth limited by N availability"

# Column B ("Category"), F ("Units"), G ("Output format") are identical for
# all ten rows.
foreach ($r in 256..265) {
    $ws.Range("B$r").Value = "livestock"
    $ws.Range("F$r").Value = "kg/y"
    $ws.Range("G$r").Value = "2f"
}

# Formatting: column E uses a centre-wrapped style on the "header" rows of
# each sub-group (milk/meat/manure N-excretion groupings), column F is
# centred + wrapped on every row, and every new row keeps the sheet's usual
# row height.
foreach ($r in 256..265) {
    $ws.Range("F$r").HorizontalAlignment = -4108
    $ws.Range("F$r").VerticalAlignment = -4108
    $ws.Range("F$r").WrapText = $true
    $ws.Rows($r).RowHeight = 20.1
}

foreach ($r in 256, 259, 262, 263) {
    $ws.Range("E$r").VerticalAlignment = -4108
    $ws.Range("E$r").WrapText = $true
}

# Column A was widened to fit the new, longer variable names.
$ws.Columns("A").ColumnWidth = 26.6

# Keep "fit to page" off (an explicit 0/0 instead of Excel's implicit
# default) as in the original page setup.
$ws.PageSetup.FitToPagesWide = 0
$ws.PageSetup.FitToPagesTall = 0

# Update the saved selection to reflect where the user finished editing.
$ws.Range("G268").Select()
